$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.315.71'
$ws.Range('E2').Value = '  -3.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.296.08'
$ws.Range('E3').Value = '  -4.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.22'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.78'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.297.25'
$ws.Range('E9').Value = '  -4.72%  '
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.48'
$ws.Range('E11').Value = '  -2.31%  '
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.331'
$ws.Range('E13').Value = '  -4.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.60'
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.711.60'
$ws.Range('E15').Value = '  -4.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.364.79'
$ws.Range('E16').Value = '  -2.94%  '
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.281.42'
$ws.Range('E18').Value = '  -6.62%  '
$ws.Range('E19').Value = '  -5.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.27'
$ws.Range('E20').Value = '  -5.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.08'
$ws.Range('E21').Value = '  -4.43%  '
$ws.Range('E22').Value = '  -4.95%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.14'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  -2.81%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.98'
$ws.Range('E27').Value = '  -7.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.29'
$ws.Range('E28').Value = '  -7.31%  '
$ws.Range('E29').Value = '  -3.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.18'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0717'
$ws.Range('E31').Value = '  -6.93%  '
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.73'
$ws.Range('E33').Value = '  -6.49%  '
$ws.Range('E34').Value = '  -5.86%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -3.72%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('E38').Value = '  -7.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.92'
$ws.Range('E39').Value = '  -6.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.08'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.49'
$ws.Range('E41').Value = '  -6.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '289.53'
$ws.Range('E42').Value = '  -10.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.98'
$ws.Range('E43').Value = '  -5.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.41'
$ws.Range('E44').Value = '  -4.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0953'
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0501'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.554'
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.22'
$ws.Range('E48').Value = '  -8.17%  '
$ws.Range('E49').Value = '  -3.87%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.00'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.52'
$ws.Range('E51').Value = '  -4.69%  '
